$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$longList = "GDG401, WMT201, DRA301, GDS301, CAA201, PST201, AFA201, GDF101, HOA101, DRF201, VCM201, DRS101, DGT101, PST201, DRP101, VNC101"

$ws.Range("A51").Value = "HB/202R"
$ws.Range("B51").Value = "Free"
$ws.Range("C51").Value = 30
$ws.Range("C51").HorizontalAlignment = -4108
$ws.Range("D51").Value = "LAB101, LAB211"
$ws.Range("E51").Value = "HB"

$ws.Range("D48").Value = $longList
$ws.Range("D50").Value = $longList

$ws.Columns("D").ColumnWidth = 126.85546875

$ws.Range("D53").Select()
